# "Generate Report for Handback"
# The handback XLIFF files for zh-cn and de-de have now come back in sync
# with en-US, so the localization-status report needs to reflect that:
#   - Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#     (Overview sheet's per-language columns, and each language sheet's Status column)
#   - Each language sheet's "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get populated with the handback info
#   - A couple of columns get widened so the longer text fits

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$targetFile = "a.md"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/643dbef8e24dfdef65392a6975d144e216fb6aee/e2e/a.md"
$bMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/643dbef8e24dfdef65392a6975d144e216fb6aee/e2e/b.md"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Widen the now-longer status columns
$wsOverview.Range("E1:F1").ColumnWidth = 29.16

# --- zh-cn sheet ---
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = $targetFile
$wsZh.Range("I3").Value = $targetFile

$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-09-01 12:41:16"
$wsZh.Range("K3").Value = "2016-09-01 12:41:16"

$wsZh.Range("C1").ColumnWidth = 29.16
$wsZh.Range("J1").ColumnWidth = 39.17

# Rebuild hyperlinks for zh-cn so A2/A3 keep their links and I2/I3 (the new
# "Latest Target File" values) become links to the source file too.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $aMdUrl, $null, $null, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $aMdUrl, $null, $null, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $bMdUrl, $null, $null, "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aMdUrl, $null, $null, "a.md")

# --- de-de sheet ---
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = $targetFile
$wsDe.Range("I3").Value = $targetFile

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("K2").Value = "2016-09-01 12:41:23"
$wsDe.Range("K3").Value = "2016-09-01 12:41:23"

$wsDe.Range("C1").ColumnWidth = 29.16
$wsDe.Range("J1").ColumnWidth = 39.17

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $aMdUrl, $null, $null, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $aMdUrl, $null, $null, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $bMdUrl, $null, $null, "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aMdUrl, $null, $null, "a.md")
